$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.870.95'
$ws.Range('E2').Value = '  +1.02%  '

$ws.Range('D3').Value = '3.316.92'
$ws.Range('E3').Value = '  +5.48%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.14%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = '3.314.88'
$ws.Range('E8').Value = '  +5.52%  '

$ws.Range('E9').Value = '  +0.86%  '

$ws.Range('E10').Value = '  +2.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.83%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.472'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.74%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000249'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.79'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.68%  '

$ws.Range('D15').Value = '3.854.63'
$ws.Range('E15').Value = '  +5.27%  '

$ws.Range('E16').Value = '  +0.18%  '

$ws.Range('D17').Value = '3.311.17'
$ws.Range('E17').Value = '  +5.37%  '

$ws.Range('D18').Value = '63.934.59'
$ws.Range('E18').Value = '  +1.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.13%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '481.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.42%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.82%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.734'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.00%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.20%  '

$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('E27').Value = '  +2.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.10%  '

$ws.Range('E29').Value = '  -0.08%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.15'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.47%  '

$ws.Range('E31').Value = '  +2.43%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.43%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.107'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.56%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.56'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.71%  '

$ws.Range('E35').Value = '  +3.65%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.67%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.11%  '

$ws.Range('D38').Value = '0.0₃0747'
$ws.Range('E38').Value = '  +6.93%  '

$ws.Range('E39').Value = '  +3.49%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '433.92'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.80%  '

$ws.Range('D41').Value = '3.075.60'
$ws.Range('E41').Value = '  +5.45%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.77'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.37'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.66%  '

$ws.Range('E44').Value = '  +2.17%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.266'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.93%  '

$ws.Range('E46').Value = '  +3.64%  '

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.08%  '

$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.09'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +13.36%  '

$ws.Range('E51').Value = '  +1.49%  '
